$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(54, 1).Value = 7143371
$ws.Cells.Item(54, 2).Value = 'NOEMI LILA LOPEZ CANEDO'
$ws.Cells.Item(54, 11).Value = 287.47

$ws.Cells.Item(55, 1).Value = 5787275
$ws.Cells.Item(55, 2).Value = 'POLONIA JEREZ RUEDA'
$ws.Cells.Item(55, 11).Value = 1010.6

$ws.Cells.Item(66, 1).Value = 7143371
$ws.Cells.Item(66, 2).Value = 'NOEMI LILA LOPEZ CANEDO'
$ws.Cells.Item(66, 11).Value = 287.47

$ws.Cells.Item(67, 1).Value = 7254095
$ws.Cells.Item(67, 2).Value = 'ROLY ALEJANDRO VELASCO VACA'
$ws.Cells.Item(67, 11).Value = 185.1

$ws.Cells.Item(80, 1).Value = 7254095
$ws.Cells.Item(80, 2).Value = 'ROLY ALEJANDRO VELASCO VACA'
$ws.Cells.Item(80, 11).Value = 185.1

$ws.Cells.Item(81, 1).Value = 7259089
$ws.Cells.Item(81, 2).Value = 'MARISOL DANIELA RAMOS'
$ws.Cells.Item(81, 11).Value = 818.38

$ws.Cells.Item(83, 1).Value = 5787275
$ws.Cells.Item(83, 2).Value = 'POLONIA JEREZ RUEDA'
$ws.Cells.Item(83, 11).Value = 1010.6

$ws.Cells.Item(84, 1).Value = 5684823
$ws.Cells.Item(84, 2).Value = 'SOLEDAD DOMINGA VASQUEZ FERNANDEZ'
$ws.Cells.Item(84, 11).Value = 665.0

$ws.Cells.Item(99, 1).Value = 7259089
$ws.Cells.Item(99, 2).Value = 'MARISOL DANIELA RAMOS'
$ws.Cells.Item(99, 11).Value = 818.38

$ws.Cells.Item(100, 1).Value = 7254095
$ws.Cells.Item(100, 2).Value = 'ROLY ALEJANDRO VELASCO VACA'
$ws.Cells.Item(100, 11).Value = 185.1

$ws.Cells.Item(102, 1).Value = 7254095
$ws.Cells.Item(102, 2).Value = 'ROLY ALEJANDRO VELASCO VACA'
$ws.Cells.Item(102, 11).Value = 185.1

$ws.Cells.Item(103, 1).Value = 5684823
$ws.Cells.Item(103, 2).Value = 'SOLEDAD DOMINGA VASQUEZ FERNANDEZ'
$ws.Cells.Item(103, 11).Value = 665.0

$ws.Cells.Item(104, 1).Value = 7259089
$ws.Cells.Item(104, 2).Value = 'MARISOL DANIELA RAMOS'
$ws.Cells.Item(104, 11).Value = 818.38

$ws.Cells.Item(108, 1).Value = 7254095
$ws.Cells.Item(108, 2).Value = 'ROLY ALEJANDRO VELASCO VACA'
$ws.Cells.Item(108, 11).Value = 185.1

$ws.Cells.Item(109, 1).Value = 7103441
$ws.Cells.Item(109, 2).Value = 'JUANA ESTHER MICHEL CRUZ'
$ws.Cells.Item(109, 11).Value = 1480.0

$ws.Cells.Item(113, 1).Value = 5797560
$ws.Cells.Item(113, 2).Value = 'LIMBER RAYNARD GARNICA MEZZA'
$ws.Cells.Item(113, 11).Value = 320.0

$ws.Cells.Item(114, 1).Value = 7190526
$ws.Cells.Item(114, 2).Value = 'MARCO ANTONIO COLODRO'
$ws.Cells.Item(114, 11).Value = 57.0

$ws.Cells.Item(115, 1).Value = 7247912
$ws.Cells.Item(115, 2).Value = 'LISELDA MILENIA ROMERO ALARCON'
$ws.Cells.Item(115, 11).Value = 8000.0

$ws.Cells.Item(121, 1).Value = 7247912
$ws.Cells.Item(121, 2).Value = 'LISELDA MILENIA ROMERO ALARCON'
$ws.Cells.Item(121, 11).Value = 8000.0

$ws.Cells.Item(122, 1).Value = 5797560
$ws.Cells.Item(122, 2).Value = 'LIMBER RAYNARD GARNICA MEZZA'
$ws.Cells.Item(122, 11).Value = 320.0

$ws.Cells.Item(123, 1).Value = 7190526
$ws.Cells.Item(123, 2).Value = 'MARCO ANTONIO COLODRO'
$ws.Cells.Item(123, 11).Value = 57.0

$ws.PageSetup.LeftHeader = "2024-02-18 01:54:35"